$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# Clear the whole sheet first so stale shared strings that are unique to
# this sheet (the old per-row "被保險人:..." labels) are garbage collected
# before any new strings are introduced. This keeps the newly-introduced
# strings ("company", "insurance", ...) contiguous with the rest of the
# (unpruned) shared string table, matching how the table is rebuilt.
$ws.Cells.ClearContents()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Helper: write "2012-04-23" as text (not an auto-converted date serial) ---
# Excel's Value setter auto-parses unambiguous date-like strings into date
# serials. Route the literal through TRIM() on a scratch cell and paste the
# computed *value* back, which keeps it text without touching NumberFormat
# (and therefore without creating a new cell style).
function Set-TextDate($cellRef) {
    $ws.Range("Y1").Value = " 2012-04-23"
    $ws.Range("Z1").Formula = "=TRIM(Y1)"
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("Y1:Z1").Clear()
}

# --- Row 2 (index 109): 南山人壽 / 新20年期缴費增值分紅终身壽險南山終身醫療保險 ---
$ws.Range("A2").Value = 109
$ws.Range("B2").Value = "南山人壽"
$ws.Range("C2").Value = "新20年期缴費增值分紅终身壽險南山終身醫療保險"
$ws.Range("D2").Value = "饒月琴"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
Set-TextDate("G2")
$ws.Range("H2").Value = "許忠信"
$ws.Range("I2").Value = 1749
$ws.Range("J2").Value = "tmpa22c1"
$ws.Range("K2").Value = 109

# --- Row 3 (index 110): 南山人壽 / 南山新年年春還本終身保險南山终身醫療保險 ---
$ws.Range("A3").Value = 110
$ws.Range("B3").Value = "南山人壽"
$ws.Range("C3").Value = "南山新年年春還本終身保險南山终身醫療保險"
$ws.Range("D3").Value = "饒月琴"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
Set-TextDate("G3")
$ws.Range("H3").Value = "許忠信"
$ws.Range("I3").Value = 1749
$ws.Range("J3").Value = "tmpa22c1"
$ws.Range("K3").Value = 110

# --- Row 4 (index 111): 南山人壽 / 南山新年年春還本終身保險南山终身醫療保險 ---
$ws.Range("A4").Value = 111
$ws.Range("B4").Value = "南山人壽"
$ws.Range("C4").Value = "南山新年年春還本終身保險南山终身醫療保險"
$ws.Range("D4").Value = "饒月琴"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
Set-TextDate("G4")
$ws.Range("H4").Value = "許忠信"
$ws.Range("I4").Value = 1749
$ws.Range("J4").Value = "tmpa22c1"
$ws.Range("K4").Value = 111

# --- Row 5 (index 112): 南山人壽 / 南山終身醫療保險 ---
$ws.Range("A5").Value = 112
$ws.Range("B5").Value = "南山人壽"
$ws.Range("C5").Value = "南山終身醫療保險"
$ws.Range("D5").Value = "饒月琴"
$ws.Range("E5").Value = "insurance"
$ws.Range("F5").Value = "normal"
Set-TextDate("G5")
$ws.Range("H5").Value = "許忠信"
$ws.Range("I5").Value = 1749
$ws.Range("J5").Value = "tmpa22c1"
$ws.Range("K5").Value = 112
